# Adds the two new worksheets introduced by the commit:
#   "Model 1 Parameters Summary" (mirrors "Model 2 Parameters Summary")
#   "Model 1_Sig"                (mirrors a district-level summary sheet,
#                                  e.g. "Summary_FuelPov")
# Both are appended at the end of the workbook, and re-use the existing
# bold+bordered header style (cell style index 3) already present in the
# workbook by copying formats from an existing, similarly-styled sheet.

$wb = $excel.ActiveWorkbook

# A sheet that already uses the "header" style (bold, thin border, center/top
# alignment) we want to reuse for the header row (B1:G1) and the hidden
# index column (A).
$styleSource = $wb.Worksheets.Item("Model 2 Parameters Summary")

# ---------------------------------------------------------------------
# Sheet: "Model 1 Parameters Summary"
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$paramSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$paramSheet.Name = "Model 1 Parameters Summary"

# Reuse the existing bold/bordered header style for the header row and the
# numeric index column.
$styleSource.Range("B1:G1").Copy()
$paramSheet.Range("B1:G1").PasteSpecial(-4122)
$styleSource.Range("A2:A3").Copy()
$paramSheet.Range("A2:A3").PasteSpecial(-4122)

$paramHeaders = @("Variable", "Mean", "STD", "Min", "Median", "Max")
for ($i = 0; $i -lt $paramHeaders.Length; $i++) {
    $paramSheet.Cells.Item(1, $i + 2).Value = $paramHeaders[$i]
}

$paramRows = @(
    @(0, "Intercept", 5.711, 1.632, 2.011, 5.722, 8.08),
    @(1, "Combined Decile", -0.217, 0.233, -1.079, -0.168, 0.485)
)

foreach ($row in $paramRows) {
    $r = [int]$row[0] + 2
    $paramSheet.Cells.Item($r, 1).Value = $row[0]
    for ($c = 0; $c -lt 6; $c++) {
        $paramSheet.Cells.Item($r, $c + 2).Value = $row[$c + 1]
    }
}

# ---------------------------------------------------------------------
# Sheet: "Model 1_Sig"
# ---------------------------------------------------------------------
$sigSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $paramSheet)
$sigSheet.Name = "Model 1_Sig"

$styleSource.Range("B1:G1").Copy()
$sigSheet.Range("B1:I1").PasteSpecial(-4122)
$styleSource.Range("A2:A3").Copy()
$sigSheet.Range("A2:A5").PasteSpecial(-4122)

$sigHeaders = @("District", "LSOA_Count", "Mean", "Median", "Minimum", "Maximum", "Total_LSOAs", "LSOA_Coverage")
for ($i = 0; $i -lt $sigHeaders.Length; $i++) {
    $sigSheet.Cells.Item(1, $i + 2).Value = $sigHeaders[$i]
}

$sigRows = @(
    @(0, "Cherwell", 38, -0.336, -0.296, -0.629, -0.259, 102, 37.255),
    @(1, "Oxford", 51, -0.621, -0.545, -1.079, -0.298, 85, 60),
    @(2, "South Oxfordshire", 32, -0.475, -0.463, -0.795, -0.344, 93, 34.409),
    @(3, "West Oxfordshire", 4, 0.423, 0.433, 0.342, 0.485, 68, 5.882)
)

foreach ($row in $sigRows) {
    $r = [int]$row[0] + 2
    $sigSheet.Cells.Item($r, 1).Value = $row[0]
    for ($c = 0; $c -lt 8; $c++) {
        $sigSheet.Cells.Item($r, $c + 2).Value = $row[$c + 1]
    }
}
